$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H28").Value = "Erratic"
$ws.Range("H29").Value = "Fast"
$ws.Range("H32").Value = "Slow"
$ws.Range("H33").Value = "Fluctuating"
$ws.Range("H30").Value = "MediumFast"
$ws.Range("H31").Value = "MediumSlow"
$ws.Range("H31").Select()
